$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Medhy -> Mathilde typo on the existing row (row 2), and its
# street / metro columns (row 2 originally referenced "15" rue + Medhy).
$ws.Range("C2").Value = "Muthilde"
$ws.Range("D2").Value = "Rue Curdinet"
$ws.Range("E2").Value = 1
$ws.Range("J2").Value = "Curdinet"

# Add the new clients (rows 3-6)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Dorand"
$ws.Range("C3").Value = "Mothilde"
$ws.Range("D3").Value = "Rue Cordinet"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 75017
$ws.Range("G3").Value = "Paris"
$ws.Range("H3").Value = 1234567899
$ws.Range("I3").Value = "Mdorand@gmail.com"
$ws.Range("J3").Value = "Cordinet"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Darand"
$ws.Range("C4").Value = "Mathilde"
$ws.Range("D4").Value = "Rue Cardinet"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 75017
$ws.Range("G4").Value = "Paris"
$ws.Range("H4").Value = 2345678901
$ws.Range("I4").Value = "Mdarand@gmail.com"
$ws.Range("J4").Value = "Cardinet"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Dirand"
$ws.Range("C5").Value = "Mithilde"
$ws.Range("D5").Value = "Rue Cirdinet"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 75017
$ws.Range("G5").Value = "Paris"
$ws.Range("H5").Value = 3456789012
$ws.Range("I5").Value = "Mdorand@gmail.com"
$ws.Range("J5").Value = "Cirdinet"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Derand"
$ws.Range("C6").Value = "Methilde"
$ws.Range("D6").Value = "Rue Cerdinet"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 75017
$ws.Range("G6").Value = "Paris"
$ws.Range("H6").Value = 4567890123
$ws.Range("I6").Value = "Mderand@gmail.com"
$ws.Range("J6").Value = "Cerdinet"

# Hyperlinks for the new email cells (mirroring I2's mailto: hyperlink),
# and apply the hyperlink style to them.
$ws.Hyperlinks.Add($ws.Range("I3"), "mailto:Mdorand@gmail.com")
$ws.Hyperlinks.Add($ws.Range("I4"), "mailto:Mdarand@gmail.com")
$ws.Hyperlinks.Add($ws.Range("I5"), "mailto:Mdorand@gmail.com")
$ws.Hyperlinks.Add($ws.Range("I6"), "mailto:Mderand@gmail.com")

$ws.Range("I3").Style = "Lien hypertexte"
$ws.Range("I4").Style = "Lien hypertexte"
$ws.Range("I5").Style = "Lien hypertexte"
$ws.Range("I6").Style = "Lien hypertexte"

# Move the active selection the way the author left it after data entry.
$ws.Range("A7").Select()
